$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: sheet restructuring
#   before: one sheet "Sprint 1" (sheetId=1) holding the backlog table
#   after : three sheets
#     "Spring 1 (M1)"  sheetId=3  (new, near-empty "Not required" placeholder)
#     "Sprint 2 (M2)"  sheetId=1  (the original backlog data, renamed)
#     "Sprint 3 (M3)"  sheetId=2  (new, empty)
# ---------------------------------------------------------------------------
$sprint2 = $wb.ActiveSheet
$sprint2.Name = "Sprint 2 (M2)"

# New sheet inserted right after Sprint 2 (M2) -> becomes "Sprint 3 (M3)"
$sprint3 = $wb.Worksheets.Add($null, $sprint2)
$sprint3.Name = "Sprint 3 (M3)"

# New sheet inserted right before Sprint 2 (M2) -> becomes "Spring 1 (M1)" (first tab)
$spring1 = $wb.Worksheets.Add($sprint2)
$spring1.Name = "Spring 1 (M1)"

# ---------------------------------------------------------------------------
# Step 2: fill in the backlog table on "Sprint 2 (M2)"
# ---------------------------------------------------------------------------
$sprint2 = $wb.Worksheets.Item("Sprint 2 (M2)")

# Header row
$sprint2.Range("A1").Value = "Tasks"
$sprint2.Range("B1").Value = "Responsible"
$sprint2.Range("C1").Value = "Status"
$sprint2.Range("D1").Value = 1
$sprint2.Range("E1").Value = 2
$sprint2.Range("F1").Value = 3
$sprint2.Range("G1").Value = 4

# Row 2 - "Set Up Git" task (column order C, B, A matches original authoring order)
$sprint2.Range("C2").Value = "Not Started"
$sprint2.Range("B2").Value = "Hunter"
$sprint2.Range("A2").Value = "Set Up Git on team's machines"
$sprint2.Range("D2").Value = 0.5

# Rows 3-7 - "Edit and commit Person N class" block
$editTasks = @(
    "Edit and commit Person 1 class",
    "Edit and commit Person 2 class",
    "Edit and commit Person 3 class",
    "Edit and commit Person 4 class",
    "Edit and commit Person 5 class"
)
for ($i = 0; $i -lt 5; $i++) {
    $r = 3 + $i
    $sprint2.Cells.Item($r, 1).Value = $editTasks[$i]
}
$people = @("Hunter", "Bhavesh", "Stephen", "Pranil", "Naman")
for ($i = 0; $i -lt 5; $i++) {
    $r = 3 + $i
    $sprint2.Cells.Item($r, 2).Value = $people[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $r = 3 + $i
    $sprint2.Cells.Item($r, 3).Value = "Not Started"
    $sprint2.Cells.Item($r, 4).Value = 2
}

# Rows 8-12 - "Create Build File N" block
$buildTasks = @(
    "Create Build File 1",
    "Create Build File 2",
    "Create Build File 3",
    "Create Build File 4",
    "Create Build File 5"
)
for ($i = 0; $i -lt 5; $i++) {
    $r = 8 + $i
    $sprint2.Cells.Item($r, 1).Value = $buildTasks[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $r = 8 + $i
    $sprint2.Cells.Item($r, 2).Value = $people[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $r = 8 + $i
    $sprint2.Cells.Item($r, 3).Value = "Not Started"
    $sprint2.Cells.Item($r, 4).Value = 1
}

Write-Output "done"
